$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the FilesTab query text in B4: drop the `File Type` and `Breed`
# RETURN columns (and their preceding lines keep their trailing comma/space).
$newQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
 MATCH (samp:sample)-->(c) 
 WHERE samp.specific_sample_pathology IN ["Undefined"]  
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@
$newQuery = $newQuery.TrimEnd("`r", "`n")

$ws.Range("B4").Value = $newQuery

# Move the active selection on the sheet from D4 to B4.
$ws.Range("B4").Select()
